# "teaches tablosu artik doluyor"
# - Instructors sheet: tighten/extend the COURSES values for the first two rows
# - Courses sheet: add the new ELEC 204 course row (capacity 45)

$wb = $excel.ActiveWorkbook

# --- Instructors sheet: update COURSES column (E) ---
$wsInstructors = $wb.Worksheets.Item("Instructors")

# Serdar Tasiran row: "COMP 131, COMP 302" -> "COMP 131,COMP 302"
$wsInstructors.Range("E2").Value2 = "COMP 131,COMP 302"

# Ozgur Baris Akan row: "ELEC 201" -> "ELEC 201,ELEC 204"
$wsInstructors.Range("E3").Value2 = "ELEC 201,ELEC 204"

# --- Courses sheet: fill in the new ELEC 204 row ---
$wsCourses = $wb.Worksheets.Item("Courses")

# Copy A4 (the "ELEC" labeled cell, same style as the other KODU entries) down to A5
$wsCourses.Range("A4").Copy($wsCourses.Range("A5")) | Out-Null

# Fill in course number (SAYISI) and capacity (KAPASITE) for the new row
$wsCourses.Range("B5").Value2 = 204
$wsCourses.Range("C5").Value2 = 45
